$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "orinaFlav" (sheet1): add the tiempo0/tiempoF sub-header block (row8),
# fill the comment/variable table (rows 9-14), widen column C, move the two
# trailing anova labels down, and change the selection.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("orinaFlav")

# New merged sub-header row (row 8), styled like row 2.
$ws1.Range("C8").Value = $ws1.Range("C9").Value
$ws1.Range("E8").Value = $ws1.Range("D9").Value
$ws1.Range("C8:F8").HorizontalAlignment = -4108
$ws1.Range("C8:D8").Merge()
$ws1.Range("E8:F8").Merge()

# Row 9 header: Comentario / Variable labels.
$ws1.Range("C9").Value = "Comentario"
$ws1.Range("D9").Value = "Variable"
$ws1.Range("E9").Value = "Variable"
$ws1.Range("F9").Value = "Comentario"

# Rows 10-14: cluster number (B), distribution comment (C), variable list (E)
# and full comment (F).
$ws1.Range("B10").Value = 1
$ws1.Range("C10").Value = "Distribución concentrada"
$ws1.Range("E10").Value = "Delta.Frec, Delta.Grasa"
$ws1.Range("F10").Value = "Concentrado, poca variabilidad, equilibrado en factores"

$ws1.Range("B11").Value = 2
$ws1.Range("C11").Value = "Distribución muy dispersa"
$ws1.Range("E11").Value = "ES, NG"
$ws1.Range("F11").Value = "Muy disperso, alta variabilidad y alta media ES y NS, Sexo/end equilibrado"
$ws1.Range("B11").RowHeight = 15

$ws1.Range("B12").Value = 3
$ws1.Range("C12").Value = "Distribución menos concentrada"
$ws1.Range("E12").Value = "Sexo H, ES, HE.G, NG, NS"
$ws1.Range("F12").Value = "Concentrado, valor medio bajo de ES, de HE.G, de NG y NS, antro normal, Casi todo hombres"

$ws1.Range("B13").Value = 4
$ws1.Range("C13").Value = "'--"
$ws1.Range("E13").Value = "Sexo M, HE.G, NG, SA"
$ws1.Range("F13").Value = "Disperso, valor muy alto y distinto de HE.G, parecido NG,  + Sa, casi todo mujeres"

$ws1.Range("B14").Value = 5
$ws1.Range("C14").Value = "'--"
$ws1.Range("E14").Value = "Sexo H, EG, HE.G, NG, peso, grasa IRCV, BPMAX, ST, "
$ws1.Range("F14").Value = "Muy disperso niveles bajos de EG, HE.G, NG, baja media de peso, baja media de grasa, baja medai delta.IRCV, baja media Bpmx, casi todo hombre, + ST"

# Move the trailing anova labels from rows 14/17 down to rows 19/22.
$ws1.Range("B19").Value = $ws1.Range("B14").Value
$ws1.Range("B19").Value = "anova-paired"
$ws1.Range("B22").Value = "anova 2/3 vias"
$ws1.Range("B17").ClearContents()

# Widen column C to fit the new comment text.
$ws1.Columns.Item(3).ColumnWidth = 35.7

# New selection for this sheet.
$ws1.Range("C10").Select()

# ---------------------------------------------------------------------------
# Sheet "orinaAnt" (sheet2): same sub-header block (row8), rework the
# header row, and fill the variable/comment columns D and F for rows10-12.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("orinaAnt")

$ws2.Range("C8").Value = $ws2.Range("B10").Value
$ws2.Range("E8").Value = $ws2.Range("B11").Value
$ws2.Range("C8:F8").HorizontalAlignment = -4108
$ws2.Range("C8:D8").Merge()
$ws2.Range("E8:F8").Merge()

$ws2.Range("B10").ClearContents()
$ws2.Range("B11").ClearContents()

$ws2.Range("C9").Value = "Variable"
$ws2.Range("D9").Value = "Comentario"
$ws2.Range("E9").Value = "Variable"
$ws2.Range("F9").Value = "Comentario"

$ws2.Range("B10").Value = 1
$ws2.Range("D10").Value = "Concentrado y definido, poca variabilidad, mayoría ST/SU, "
$ws2.Range("F10").Value = "Concentrado y definido"

$ws2.Range("B11").Value = 2
$ws2.Range("D11").Value = "Concentrado y definido"
$ws2.Range("F11").Value = "Definido algo disperso"

$ws2.Range("B12").Value = 3
$ws2.Range("D12").Value = "Disperso"
$ws2.Range("F12").Value = "Muy disperso"

# New selection for this sheet, and make it the active tab.
$ws2.Range("D10").Select()
$ws2.Activate()
